$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '59.312.28'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +2.55%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.181.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +1.68%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '532.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -0.14%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '141.66'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +1.92%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +0.08%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +11.15%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '7.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.61%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.440'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +7.07%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +4.26%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '3.734.48'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +1.84%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '0.139'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +1.63%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '25.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -0.05%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '0.0000171'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  +4.21%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '59.362.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +2.47%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '3.180.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +1.93%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '6.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +2.55%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '13.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +2.33%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '8.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +1.27%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '375.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +2.22%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +0.12%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +5.08%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '69.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.86%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '0.168'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -0.17%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.04%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '8.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +15.14%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '0.0₃0873'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +0.95%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '22.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +4.52%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '1.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +0.76%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '6.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -0.71%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '5.26'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +1.83%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '1.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -0.83%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '6.35'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +4.34%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '157.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -1.45%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +3.54%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = "'" + 'Hedera'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'" + '0.0716'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +6.66%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'" + 'EnergySwap'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'" + '25.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -0.07%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '2.707.93'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '1.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +1.52%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '4.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +4.35%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.727'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +3.85%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +8.46%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '39.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +3.61%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +0.06%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '3.224.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +1.69%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'" + 'Stellar'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.102'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +12.15%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Value = "'" + 'ONDO'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'" + 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'" + '0.989'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +0.87%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '6.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +1.27%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '20.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +3.08%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.762'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +2.11%  '
$ws.Range('E51').Style = 'Normal'
